$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1)
#    Columns A-E are brand-new plain headers (no explicit style => default).
#    Columns F-K keep the "(m3/s)" text (already in the shared-string table)
#    plus five new headers, all styled with the bold-ish header style that we
#    build below (fontId=1, numFmtId=0, applyFont only).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# ---------------------------------------------------------------------------
# 2. Data rows (2-10) - shifted up one row compared to the source file, and
#    with a brand-new leading "idx"/"idx2" pair. Styles s="2" (integer),
#    s="1" (text) and s="3" (2-decimal) already exist in the workbook and are
#    reused automatically by the engine, matching the original per-column
#    formatting.
# ---------------------------------------------------------------------------

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 304600
$ws.Range("C2").Value = "Perlen 1 (WTA HF, in Buchrain)"
$ws.Range("D2").Value = 1873
$ws.Range("E2").Value = 1981
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.96
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 8

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 304700
$ws.Range("C3").Value = "Perlen 2 (WTA PF, in Root)"
$ws.Range("D3").Value = 1875
$ws.Range("E3").Value = 2000
$ws.Range("F3").Value = 45
$ws.Range("G3").Value = 1.15
$ws.Range("H3").Value = 1.09
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 3.9
$ws.Range("K3").Value = 7.8

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 304300
$ws.Range("C4").Value = "Thorenberg"
$ws.Range("D4").Value = 1886
$ws.Range("E4").Value = 2000
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 0.79
$ws.Range("H4").Value = 0.75
$ws.Range("I4").Value = 2.2
$ws.Range("J4").Value = 2.6
$ws.Range("K4").Value = 4.8

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 304500
$ws.Range("C5").Value = "Rathausen"
$ws.Range("D5").Value = 1896
$ws.Range("E5").Value = 1980
$ws.Range("F5").Value = 45
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 7.88
$ws.Range("J5").Value = 8.02
$ws.Range("K5").Value = 15.9

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 304200
$ws.Range("C6").Value = "Wolhusen (Geistlich)"
$ws.Range("D6").Value = 1906
$ws.Range("E6").Value = 2003
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 0.33
$ws.Range("H6").Value = 0.3
$ws.Range("I6").Value = 0.53
$ws.Range("J6").Value = 0.63
$ws.Range("K6").Value = 1.16

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 304400
$ws.Range("C7").Value = "Emmenweid"
$ws.Range("D7").Value = 1931
$ws.Range("E7").Value = 2003
$ws.Range("F7").Value = 12
$ws.Range("G7").Value = 1.08
$ws.Range("H7").Value = 1.08
$ws.Range("I7").Value = 1.59
$ws.Range("J7").Value = 2.96
$ws.Range("K7").Value = 4.55

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 304150
$ws.Range("C8").Value = "Mühlenplatz"
$ws.Range("D8").Value = 1998
$ws.Range("E8").ClearContents()
$ws.Range("F8").Value = 58
$ws.Range("G8").Value = 0.93
$ws.Range("H8").Value = 0.83
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 3

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 304250
$ws.Range("C9").Value = "Stollen"
$ws.Range("D9").Value = 1999
$ws.Range("E9").ClearContents()
$ws.Range("F9").Value = 0.35
$ws.Range("G9").Value = 0.71
$ws.Range("H9").Value = 0.64
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 3

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 304225
$ws.Range("C10").Value = "Ettisbühl"
$ws.Range("D10").Value = 2011
$ws.Range("E10").ClearContents()
$ws.Range("F10").Value = 16
$ws.Range("G10").Value = 0.87
$ws.Range("H10").Value = 0.87
$ws.Range("I10").Value = 1.72
$ws.Range("J10").Value = 2.78
$ws.Range("K10").Value = 4.5

# ---------------------------------------------------------------------------
# 3. Remove what used to be row 11 (the sheet now only spans down to row 10).
# ---------------------------------------------------------------------------
$ws.Rows.Item(11).Delete()

# ---------------------------------------------------------------------------
# 4. Header style for F1:K1 - a cellXf with fontId=1 (existing 9pt Arial font)
#    and applyFont only (no applyNumberFormat). We get this exact combination
#    by creating a throw-away named style, tweaking its font to match the
#    existing 9pt Arial font, applying it, then deleting the named style -
#    which leaves the new cellXf behind (xfId falls back to 0) without adding
#    any permanent entries to cellStyleXfs/cellStyles.
# ---------------------------------------------------------------------------
$headerStyle = $wb.Styles.Add("__TmpHeaderStyle")
$headerStyle.Font.Name = "Arial"
$headerStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "__TmpHeaderStyle"
$headerStyle.Delete()

# ---------------------------------------------------------------------------
# 5. Selection / view bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("A2:K2").Select()

Write-Output "done"
